$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric cell values in rows 2-4 (capital structure database refresh)

# Row 2
$ws.Range("D2").Value = 0.09849999999999999
$ws.Range("E2").Value = 0.136
$ws.Range("F2").Value = 0.07719999999999999
$ws.Range("G2").Value = 0.3940333911139494
$ws.Range("H2").Value = 0.3827205546939147
$ws.Range("I2").Value = 0.4744092692272602
$ws.Range("J2").Value = 0.3212594559008717
$ws.Range("K2").Value = 365.5
$ws.Range("L2").Value = 0.3334549767356992
$ws.Range("M2").Value = 152.49
$ws.Range("N2").Value = 0.01908080782802372
$ws.Range("O2").Value = 0.4172093023255814
$ws.Range("P2").Value = 133.99
$ws.Range("Q2").Value = 0.01676593508346055
$ws.Range("R2").Value = 0.366593707250342
$ws.Range("S2").Value = 18.5
$ws.Range("T2").Value = 0.1213194307823464
$ws.Range("U2").Value = 750.4000000000001
$ws.Range("V2").Value = 0.09389624364974099
$ws.Range("W2").Value = 0.2754515703040688
$ws.Range("X2").Value = 0.02047265216076195
$ws.Range("Y2").Value = 0.2549789181433069
$ws.Range("Z2").Value = 0.5137410859289501
$ws.Range("AA2").Value = 0.1042502428106921
$ws.Range("AB2").Value = 0.02031053880824775
$ws.Range("AC2").Value = 0.08393970400244438
$ws.Range("AD2").Value = 1875.9
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1875.9
$ws.Range("AG2").Value = 1125.5
$ws.Range("AH2").Value = 0.1901050903452679
$ws.Range("AI2").Value = 0.5952781391806556
$ws.Range("AJ2").Value = 0.1234466344202779
$ws.Range("AK2").Value = 0.4687825398808781
$ws.Range("AL2").Value = 19.1
$ws.Range("AM2").Value = 11.84
$ws.Range("AN2").Value = 3.460431654676259
$ws.Range("AO2").Value = 27.22513089005236
$ws.Range("AP2").Value = 2.076185205681608
$ws.Range("AQ2").Value = 43.91891891891891

# Row 3
$ws.Range("D3").Value = 0.09849999999999999
$ws.Range("E3").Value = 0.136
$ws.Range("F3").Value = 0.07719999999999999
$ws.Range("G3").Value = 0.4388335704125177
$ws.Range("H3").Value = 0.4262345051818736
$ws.Range("I3").Value = 0.5283478967689493
$ws.Range("J3").Value = 0.3836349618083363
$ws.Range("K3").Value = 349.1
$ws.Range("L3").Value = 0.3547043283885389
$ws.Range("M3").Value = 148.4
$ws.Range("N3").Value = 0.01934407424787528
$ws.Range("O3").Value = 0.4250930965339444
$ws.Range("P3").Value = 129.9
$ws.Range("Q3").Value = 0.01693258251212263
$ws.Range("R3").Value = 0.3720996849040389
$ws.Range("S3").Value = 18.5
$ws.Range("T3").Value = 0.1246630727762803
$ws.Range("U3").Value = 666.2
$ws.Range("V3").Value = 0.08683977266802231
$ws.Range("W3").Value = 0.3698876880695063
$ws.Range("X3").Value = 0.01823592418503987
$ws.Range("Y3").Value = 0.3516517638844664
$ws.Range("Z3").Value = 0.5434866640896793
$ws.Range("AA3").Value = 0.2085004856213843
$ws.Range("AB3").Value = 0.01802237842405653
$ws.Range("AC3").Value = 0.1904781071973277
$ws.Range("AD3").Value = 1564.2
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 1564.2
$ws.Range("AG3").Value = 898
$ws.Range("AH3").Value = 0.1693626973299552
$ws.Range("AI3").Value = 0.5718150246755621
$ws.Range("AJ3").Value = 0.1047890216579537
$ws.Range("AK3").Value = 0.4339631759532209
$ws.Range("AL3").Value = 19.1
$ws.Range("AM3").Value = 11.84
$ws.Range("AN3").Value = 2.885445489762037
$ws.Range("AO3").Value = 27.22513089005236
$ws.Range("AP3").Value = 1.656520937096477
$ws.Range("AQ3").Value = 43.91891891891891

# Row 4
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 16.4
$ws.Range("L4").Value = 0.1465594280607685
$ws.Range("M4").Value = 4.09
$ws.Range("N4").Value = 0.01277326670830731
$ws.Range("O4").Value = 0.249390243902439
$ws.Range("P4").Value = 4.09
$ws.Range("Q4").Value = 0.01277326670830731
$ws.Range("R4").Value = 0.249390243902439
$ws.Range("U4").Value = 84.2
$ws.Range("V4").Value = 0.2629606495940038
$ws.Range("W4").Value = 0.1810154525386314
$ws.Range("X4").Value = 0.02270938013648402
$ws.Range("Y4").Value = 0.1583060724021473
$ws.Range("Z4").Value = 0.3467993119799173
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.02259869919243897
$ws.Range("AC4").Value = -0.02259869919243897
$ws.Range("AD4").Value = 311.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 311.7
$ws.Range("AG4").Value = 227.5
$ws.Range("AH4").Value = 0.4932742522551036
$ws.Range("AI4").Value = 0.7496392496392497
$ws.Range("AJ4").Value = 0.4153733795873653
$ws.Range("AK4").Value = 0.6860675512665863

# Cells removed entirely in row 4 (no longer have a value)
$ws.Range("F4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()
